$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B:K shift right to C:L)
$ws.Columns.Item(2).Insert()

# Inserting a column copies the left neighbor's (column A's bold/border)
# formatting into the new column's data cells - the data rows (2-5)
# should stay unstyled like the rest of the numeric columns.
$ws.Range("B2:B5").ClearFormats()

# The newly inserted column B header (row 1) needs the same bold/border
# style ("s=1") that the other header cells already carry. Copy the
# formatting from the neighboring header cell, then set its value.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row (row 1) text
$ws.Range("B1").Value = "Vdot_1"
$ws.Range("C1").Value = "Vdot el. Comp. in [l/s]"
$ws.Range("D1").Value = "Vdot Intercooler in [l/s]"
$ws.Range("E1").Value = "LV DCDC T_in"
$ws.Range("F1").Value = "HV DCDC T_in"
$ws.Range("G1").Value = "Inverter T_in"
$ws.Range("H1").Value = "HPDU T_in"
$ws.Range("I1").Value = "Compressor T_in"
$ws.Range("J1").Value = "Intercooler T_in"
$ws.Range("K1").Value = "BoP T_out"
$ws.Range("L1").Value = "P_out in [bar]"

# Row 2 data
$ws.Range("B2").Value = 0.25
$ws.Range("C2").Value = 0.08059539500474169
$ws.Range("D2").Value = 0.1694046049952583
$ws.Range("E2").Value = 312.1704476099254
$ws.Range("F2").Value = 312.3507121569987
$ws.Range("G2").Value = 313.4332218419817
$ws.Range("H2").Value = 312.1704476099254
$ws.Range("I2").Value = 315.0110017603065
$ws.Range("J2").Value = 313.3119395659741
$ws.Range("K2").Value = 324.836433382633
$ws.Range("L2").Value = 0.7506248225323764

# Row 3 data
$ws.Range("B3").Value = 0.3333333333333333
$ws.Range("C3").Value = 0.1059956085962485
$ws.Range("D3").Value = 0.2273373914037499
$ws.Range("E3").Value = 314.9177038178762
$ws.Range("F3").Value = 315.0519199006713
$ws.Range("G3").Value = 315.877150067895
$ws.Range("H3").Value = 314.9177038178762
$ws.Range("I3").Value = 317.0761544598553
$ws.Range("J3").Value = 315.7676611761784
$ws.Range("K3").Value = 324.4148216347054
$ws.Range("L3").Value = 0.5915138621357861

# Row 4 data
$ws.Range("B4").Value = 0.4166666666666667
$ws.Range("C4").Value = 0.1313263722084737
$ws.Range("D4").Value = 0.2853406277915264
$ws.Range("E4").Value = 316.5650743986903
$ws.Range("F4").Value = 316.6719656471541
$ws.Range("G4").Value = 317.3391841463729
$ws.Range("H4").Value = 316.5650743986903
$ws.Range("I4").Value = 318.3066606681732
$ws.Range("J4").Value = 317.2420121325483
$ws.Range("K4").Value = 324.1618541532335
$ws.Range("L4").Value = 0.394463876324013

# Row 5 data
$ws.Range("B5").Value = 0.5
$ws.Range("C5").Value = 0.1566188491071225
$ws.Range("D5").Value = 0.3433811508928775
$ws.Range("E5").Value = 317.6629731430928
$ws.Range("F5").Value = 317.7517777873888
$ws.Range("G5").Value = 318.3119444875464
$ws.Range("H5").Value = 317.6629731430928
$ws.Range("I5").Value = 319.1230630964595
$ws.Range("J5").Value = 318.2253790890475
$ws.Range("K5").Value = 323.9932120428591
$ws.Range("L5").Value = 0.1594804941517815
